$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new reviewer row data: email hyperlink cell D8 ("armonravid@gmail.com"),
# mirroring the existing email cells (C/D columns) in the sheet.
$ws.Range("D8").Value = "armonravid@gmail.com"
$ws.Hyperlinks.Add($ws.Range("D8"), "mailto:armonravid@gmail.com", "", "", "armonravid@gmail.com")

# Match the formatting used by the other email cells in column D (e.g. D7)
# instead of Excel's default auto-applied "Hyperlink" style.
$ws.Range("D8").Font.Name = $ws.Range("D7").Font.Name
$ws.Range("D8").Font.Size = $ws.Range("D7").Font.Size
$ws.Range("D8").Font.Color = $ws.Range("D7").Font.Color
$ws.Range("D8").Font.Underline = $ws.Range("D7").Font.Underline

# Move the active selection to D9, as recorded in the saved workbook view.
$ws.Range("D9").Select()
